# Insert a brand-new record row before the existing row 50.
# This shifts the previous rows 50..80 down to 51..81 (preserving all of
# their data/styles), matching the target diff which effectively pushes
# every existing entry down by one row and prepends a new weekly entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("50").Insert()

# Populate the newly inserted row 50 with the new weekly price entry.
$ws.Cells.Item(50, 1).Value() = 7
$ws.Cells.Item(50, 2).Value() = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(50, 3).Value() = "Ñuble"
$ws.Cells.Item(50, 4).Value() = 44818
$ws.Cells.Item(50, 5).Value() = 16
$ws.Cells.Item(50, 6).Value() = "Fruta"
$ws.Cells.Item(50, 7).Value() = 100108
$ws.Cells.Item(50, 8).Value() = "Tropicales y subtropicales"
$ws.Cells.Item(50, 9).Value() = 100108002
$ws.Cells.Item(50, 10).Value() = "Mango"
$ws.Cells.Item(50, 11).Value() = "Sin especificar"
$ws.Cells.Item(50, 12).Value() = "Primera"
$ws.Cells.Item(50, 13).Value() = 120
$ws.Cells.Item(50, 14).Value() = 9000
$ws.Cells.Item(50, 15).Value() = 10000
$ws.Cells.Item(50, 16).Value() = 9500
$ws.Cells.Item(50, 17).Value() = "`$/bandeja 4 kilos"
$ws.Cells.Item(50, 18).Value() = "Brasil"
$ws.Cells.Item(50, 19).Value() = 2375
$ws.Cells.Item(50, 20).Value() = 4
